$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Workbook-level summary: Total Students count (K4/L4 block) increments by one
#    as a new student joined group B1E1.
$ws.Range("L4").Value = 321

# 2) Group B1E1 roster grows from 25 to 26 students: the per-session attendance
#    fractions ("x/25" -> "x/26") are recalculated for every one of its 27 sessions,
#    plus the group summary row (M23/S23) and a few "Recorded By" cells that pick up
#    the "System" auto-recorder alongside the existing grader once the roster changed.
$ws.Range("H212").Value = "20/26"
$ws.Range("H213").Value = "19/26"
$ws.Range("H214").Value = "3/26"
$ws.Range("H215").Value = "24/26"
$ws.Range("H216").Value = "20/26"
$ws.Range("H217").Value = "24/26"
$ws.Range("G218").Value = "dnasr281@gmail.com, System"
$ws.Range("H218").Value = "26/26"
$ws.Range("H219").Value = "24/26"
$ws.Range("H220").Value = "25/26"
$ws.Range("G221").Value = "dnasr281@gmail.com, System"
$ws.Range("H221").Value = "25/26"
$ws.Range("H222").Value = "21/26"
$ws.Range("H223").Value = "22/26"
$ws.Range("G224").Value = "dnasr281@gmail.com, System"
$ws.Range("H224").Value = "22/26"
$ws.Range("H225").Value = "20/26"
$ws.Range("H226").Value = "18/26"
$ws.Range("H227").Value = "0/26"
$ws.Range("H228").Value = "0/26"
$ws.Range("H229").Value = "0/26"
$ws.Range("H230").Value = "0/26"
$ws.Range("H231").Value = "0/26"
$ws.Range("H232").Value = "0/26"
$ws.Range("H233").Value = "0/26"
$ws.Range("H234").Value = "0/26"
$ws.Range("H235").Value = "0/26"
$ws.Range("H236").Value = "0/26"
$ws.Range("H237").Value = "0/26"
$ws.Range("H238").Value = "0/26"

# Group-level rollup for B1E1 (row 23): Students 25 -> 26, Avg Attendance % refreshed
$ws.Range("M23").Value = 26
# S23 holds a text percentage ("80.3%"); a leading apostrophe keeps it text instead of
# Excel auto-converting it to a numeric percent, then we re-apply R23's format (same
# style already shared by the whole row) so the cell's look is untouched.
$ws.Range("S23").Value = "'80.3%"
$ws.Range("R23").Copy()
$ws.Range("S23").PasteSpecial(-4122)

# 3) Sheet-wide normalization: "Recorded By" values listing both the System auto-recorder
#    and the human grader are re-ordered/re-joined consistently as "<email>, System"
#    across every affected session row in every group.
$ws.Range("G8").Value = "dnasr281@gmail.com, System"
$ws.Range("G9").Value = "dnasr281@gmail.com, System"
$ws.Range("G10").Value = "dnasr281@gmail.com, System"
$ws.Range("G12").Value = "dnasr281@gmail.com, System"
$ws.Range("G14").Value = "dnasr281@gmail.com, System"
$ws.Range("G15").Value = "dnasr281@gmail.com, System"
$ws.Range("G17").Value = "dnasr281@gmail.com, System"
$ws.Range("G34").Value = "dnasr281@gmail.com, System"
$ws.Range("G35").Value = "dnasr281@gmail.com, System"
$ws.Range("G36").Value = "dnasr281@gmail.com, System"
$ws.Range("G38").Value = "dnasr281@gmail.com, System"
$ws.Range("G40").Value = "dnasr281@gmail.com, System"
$ws.Range("G41").Value = "dnasr281@gmail.com, System"
$ws.Range("G43").Value = "dnasr281@gmail.com, System"
$ws.Range("G60").Value = "dnasr281@gmail.com, System"
$ws.Range("G61").Value = "dnasr281@gmail.com, System"
$ws.Range("G62").Value = "dnasr281@gmail.com, System"
$ws.Range("G64").Value = "dnasr281@gmail.com, System"
$ws.Range("G66").Value = "dnasr281@gmail.com, System"
$ws.Range("G67").Value = "dnasr281@gmail.com, System"
$ws.Range("G69").Value = "dnasr281@gmail.com, System"
$ws.Range("G86").Value = "dnasr281@gmail.com, System"
$ws.Range("G87").Value = "dnasr281@gmail.com, System"
$ws.Range("G88").Value = "dnasr281@gmail.com, System"
$ws.Range("G90").Value = "dnasr281@gmail.com, System"
$ws.Range("G92").Value = "dnasr281@gmail.com, System"
$ws.Range("G93").Value = "dnasr281@gmail.com, System"
$ws.Range("G95").Value = "dnasr281@gmail.com, System"
$ws.Range("G112").Value = "dnasr281@gmail.com, System"
$ws.Range("G113").Value = "dnasr281@gmail.com, System"
$ws.Range("G114").Value = "dnasr281@gmail.com, System"
$ws.Range("G116").Value = "dnasr281@gmail.com, System"
$ws.Range("G118").Value = "dnasr281@gmail.com, System"
$ws.Range("G119").Value = "dnasr281@gmail.com, System"
$ws.Range("G121").Value = "dnasr281@gmail.com, System"
$ws.Range("G138").Value = "dnasr281@gmail.com, System"
$ws.Range("G139").Value = "dnasr281@gmail.com, System"
$ws.Range("G140").Value = "dnasr281@gmail.com, System"
$ws.Range("G142").Value = "dnasr281@gmail.com, System"
$ws.Range("G144").Value = "dnasr281@gmail.com, System"
$ws.Range("G145").Value = "dnasr281@gmail.com, System"
$ws.Range("G147").Value = "dnasr281@gmail.com, System"
$ws.Range("G164").Value = "dnasr281@gmail.com, System"
$ws.Range("G167").Value = "dnasr281@gmail.com, System"
$ws.Range("G170").Value = "dnasr281@gmail.com, System"
$ws.Range("G191").Value = "dnasr281@gmail.com, System"
$ws.Range("G194").Value = "dnasr281@gmail.com, System"
$ws.Range("G197").Value = "dnasr281@gmail.com, System"
$ws.Range("G245").Value = "dnasr281@gmail.com, System"
$ws.Range("G248").Value = "dnasr281@gmail.com, System"
$ws.Range("G251").Value = "dnasr281@gmail.com, System"
$ws.Range("G272").Value = "dnasr281@gmail.com, System"
$ws.Range("G275").Value = "dnasr281@gmail.com, System"
$ws.Range("G278").Value = "dnasr281@gmail.com, System"
$ws.Range("G299").Value = "dnasr281@gmail.com, System"
$ws.Range("G302").Value = "dnasr281@gmail.com, System"
$ws.Range("G305").Value = "dnasr281@gmail.com, System"
